# Update the "Marking" row's Right value (B11): 3 -> 5
# Update the "Total" row's Right value (B12): 60 -> 100
# Update the "Total" row's Max "correct/total" text (E12): "56/84" -> "100/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100/140"
